$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper "templates" borrowed from existing content already in the
# document so newly-created runs get exactly the right rPr/pPr
# (the COM layer merges/inherits formatting from neighboring text
# when using plain InsertAfter, so we copy FormattedText from a
# suitable donor range instead, then overwrite the text in place).
# ------------------------------------------------------------------

# Run with NO rPr at all (same shape as the existing "survival guide" run).
$pGuide = $d.Paragraphs(12)
$noLangTemplate = $d.Range($pGuide.Range.Start + 8, $pGuide.Range.End - 1)

# Run/paragraph with rPr lang=bg-BG (single run paragraph "Име:").
$pName = $d.Paragraphs(25)
$bgTemplate = $d.Range($pName.Range.Start, $pName.Range.End - 1)

# Paragraph formatting template: ListParagraph + numPr(ilvl0,numId12) + ind left=1276 + rPr lang=bg-BG
$pBeer = $d.Paragraphs(32)
$pPrTemplateRange = $d.Range($pBeer.Range.Start, $pBeer.Range.End)

# ------------------------------------------------------------------
# 1) Paragraph "И понеже е GUIDE ..." -> append " (Vector)"
#    " (" and ")" keep bg-BG language, "Vector" gets no rPr at all.
# ------------------------------------------------------------------
$pGuideLine = $d.Paragraphs(33)
$insertPoint = $pGuideLine.Range.End - 1

$bgLen = $bgTemplate.End - $bgTemplate.Start

$rOpen = $d.Range($insertPoint, $insertPoint)
$rOpen.FormattedText = $bgTemplate.FormattedText
$rOpenFixed = $d.Range($insertPoint, $insertPoint + $bgLen)
$rOpenFixed.Text = " ("

$vecStart = $rOpenFixed.End
$noLangLen = $noLangTemplate.End - $noLangTemplate.Start
$rVec = $d.Range($vecStart, $vecStart)
$rVec.FormattedText = $noLangTemplate.FormattedText
$rVecFixed = $d.Range($vecStart, $vecStart + $noLangLen)
$rVecFixed.Text = "Vector"

$closeStart = $rVecFixed.End
$rClose = $d.Range($closeStart, $closeStart)
$rClose.FormattedText = $bgTemplate.FormattedText
$rCloseFixed = $d.Range($closeStart, $closeStart + $bgLen)
$rCloseFixed.Text = ")"

# ------------------------------------------------------------------
# 2) Replace the old "5.Интервюта със студенти" paragraph with two new
#    bullet paragraphs (numId=12 list), then delete the old paragraph.
# ------------------------------------------------------------------
$pOld = $d.Paragraphs(34)
$insertAt = $pOld.Range.Start

# Clone #1 (formatting donor) -> becomes "Лого на сайта (Vector)"
$clone1 = $d.Range($insertAt, $insertAt)
$clone1.FormattedText = $pPrTemplateRange.FormattedText
$pNew1 = $d.Paragraphs(34)

# Clone #2 (formatting donor) -> becomes "Схеми от рода на ..."
$insertAt2 = $pNew1.Range.End
$clone2 = $d.Range($insertAt2, $insertAt2)
$clone2.FormattedText = $pPrTemplateRange.FormattedText
$pNew2 = $d.Paragraphs(35)

# --- Fill paragraph 1: "Лого на сайта (" + "Vector)" ---
$p1Start = $pNew1.Range.Start
$p1TextLen = $pNew1.Range.End - $pNew1.Range.Start - 1
$p1Rng = $d.Range($p1Start, $p1Start + $p1TextLen)
$p1Rng.Text = "Лого на сайта ("

$vec2Start = $p1Rng.End
$rVec2 = $d.Range($vec2Start, $vec2Start)
$rVec2.FormattedText = $noLangTemplate.FormattedText
$rVec2Fixed = $d.Range($vec2Start, $vec2Start + $noLangLen)
$rVec2Fixed.Text = "Vector)"

# --- Fill paragraph 2: "Схеми от рода на бира + бира = много код, бира + бира на втора = малко код." ---
$pNew2 = $d.Paragraphs(35)
$p2Start = $pNew2.Range.Start
$p2TextLen = $pNew2.Range.End - $pNew2.Range.Start - 1
$p2Rng = $d.Range($p2Start, $p2Start + $p2TextLen)
$p2Rng.Text = "Схеми от рода на бира + бира = много код, бира + бира на втора = малко код."

# ------------------------------------------------------------------
# 3) Move the _GoBack bookmark to the end of the new "Схеми ..."
#    paragraph (this also removes the old _GoBack near the top of the
#    document, since _GoBack is a singleton bookmark in Word).
#    Bookmarks.Add placed on a zero-length range is unreliable in this
#    host, so we bookmark a throw-away placeholder character and then
#    delete just that character, leaving the (now zero-width) bookmark
#    correctly anchored in place.
# ------------------------------------------------------------------
$bmAnchor = $p2Rng.End
$placeholder = $d.Range($bmAnchor, $bmAnchor)
$placeholder.InsertAfter("X")
$placeholderFixed = $d.Range($bmAnchor, $bmAnchor + 1)
$d.Bookmarks.Add("_GoBack", $placeholderFixed)
$bm = $d.Bookmarks("_GoBack")
$bmRange = $d.Range($bm.Start, $bm.End)
$bmRange.Delete()

# ------------------------------------------------------------------
# 4) Remove the old "5.Интервюта със студенти" paragraph (now pushed
#    down to index 36).
# ------------------------------------------------------------------
$pOldNow = $d.Paragraphs(36)
$pOldNow.Range.Delete()

Write-Output "done"
